$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 178: "Nb nouveaux cas positifs" corrected from 17 to 18
$ws.Range("C178").Value = 18

# Row 181: "Nb nouveaux cas positifs" corrected from 15 to 14
$ws.Range("C181").Value = 14

# Row 182: "Nb nouveaux cas positifs" corrected from 10 to 13
$ws.Range("C182").Value = 13

# Row 183: "Nb nouveaux cas positifs" filled in as 9 (was blank)
$ws.Range("C183").Value = 9

# Row 184 (2020-08-27): fill in the day's figures (was fully blank)
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 1
$ws.Range("F184").Value = 1
$ws.Range("G184").Value = 3
$ws.Range("I184").Value = 0
$ws.Range("L184").Value = "0"
$ws.Range("M184").Value = "0"

# Move the active selection on the frozen pane to reflect where entry continued
$ws.Range("O186").Select() | Out-Null
